$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 5, 8, 9
$ws.Range("F2").Value = -1
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -2
